$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column G ---
$ws.Columns("G").ColumnWidth = 21.14

# --- New "ICR1 (H/L)" block header ---
$ws.Range("K2").Value = "ICR1 (H/L)"

# --- New column headers (right aligned) ---
$ws.Range("K3").Value = "top start"
$ws.Range("K3").HorizontalAlignment = -4152   # xlRight

$ws.Range("O3").Value = "top end"
$ws.Range("O3").HorizontalAlignment = -4152   # xlRight

$ws.Range("L3").Value = "duty"
$ws.Range("L3").HorizontalAlignment = -4152   # xlRight

$ws.Range("M3").Value = "freq"
$ws.Range("M3").HorizontalAlignment = -4152   # xlRight

# Q3 uses text-quote-prefixed "freq" with right alignment
$ws.Range("Q3").Value = "'freq"
$ws.Range("Q3").HorizontalAlignment = -4152   # xlRight

# --- Column G header relabel: "time unit" -> "top time" (added last) ---
$ws.Range("G3").Value = "top time"

# --- New data rows ---
$ws.Range("K5").Value = 16
$ws.Range("L5").Value = 16
$ws.Range("M5").Value = 1000000
$ws.Range("O5").Value = 65536
$ws.Range("Q5").Value = 244

$ws.Range("K7").Value = 8163
$ws.Range("M7").Value = 244
$ws.Range("O7").Value = 65536
$ws.Range("Q7").Value = 30

# --- D9 changed from 65536 to 1 ---
$ws.Range("D9").Value = 1

$ws.Range("K9").Value = 1041
$ws.Range("M9").Value = 30
$ws.Range("O9").Value = 31250
$ws.Range("Q9").Value = 1

# --- Row 13 new duty-cycle helper cells ---
$ws.Range("J13").Value = 50
$ws.Range("K13").Formula = "=1/J13"
$ws.Range("L13").Formula = "=K13/G9"

# --- Selection / active cell moves to L13 ---
$ws.Range("L13").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
